$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "r876"
$ws.Range("B7").Value = "vicky"
$ws.Range("C7").Value = "are we back to normal?"
$ws.Range("D7").Value = "2025-09-30 20:30:57"
